# Insert a new weekly price record at row 150 for "Vega Modelo de Temuco" /
# Coliflor. This shifts the existing rows 150-223 down to 151-224 (the old
# row 223 becomes the new row 224), and the freshly inserted row 150 is
# populated with a new observation (same market/category/etc. as the row
# that used to occupy position 150, but with its own date and volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 150, pushing everything
# below (through row 223) down by one row.
$ws.Rows.Item(150).EntireRow.Insert()

# Populate the newly inserted row 150 with the new record.
$ws.Cells.Item(150, 1).Value = 10
$ws.Cells.Item(150, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(150, 3).Value = "La Araucanía"
$ws.Cells.Item(150, 4).Value = 44455
$ws.Cells.Item(150, 5).Value = 9
$ws.Cells.Item(150, 6).Value = 100112008
$ws.Cells.Item(150, 7).Value = "Coliflor"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 900
$ws.Cells.Item(150, 11).Value = 800
$ws.Cells.Item(150, 12).Value = 800
$ws.Cells.Item(150, 13).Value = 800
$ws.Cells.Item(150, 14).Value = "`$/unidad"
$ws.Cells.Item(150, 15).Value = "Región Metropolitana"
$ws.Cells.Item(150, 16).Value = 800
$ws.Cells.Item(150, 17).Value = 1
$ws.Cells.Item(150, 18).Value = "Hortaliza"
